$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.5364909999999999
$ws.Range("H2").Value = 1.609473
$ws.Range("M2").Value = 0.177232
$ws.Range("N2").Value = 0.5316959999999999
$ws.Range("O2").Value = 0.0005104719838156216
$ws.Range("P2").Value = 0.0005104719838156217
$ws.Range("Q2").Value = 0.09508337291199997
$ws.Range("R2").Value = 0.8557503562079999
$ws.Range("S2").Value = 0.0005104719838156216
$ws.Range("T2").Value = 0.0005104719838156217

# Row 3
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.5364909999999999
$ws.Range("H3").Value = 1.609473
$ws.Range("O3").Value = 0.0002336907822601807
$ws.Range("P3").Value = 0.0002336907822601807
$ws.Range("Q3").Value = 0.04352855494566667
$ws.Range("R3").Value = 0.391756994511
$ws.Range("S3").Value = 0.0002336907822601807
$ws.Range("T3").Value = 0.0002336907822601807

# Row 4
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.5364909999999999
$ws.Range("H4").Value = 1.609473
$ws.Range("M4").Value = 274.5137023333334
$ws.Range("N4").Value = 823.541107
$ws.Range("O4").Value = 0.7906673411949746
$ws.Range("P4").Value = 0.7906673411949746
$ws.Range("Q4").Value = 147.2741306785123
$ws.Range("R4").Value = 1325.467176106611
$ws.Range("S4").Value = 0.7906673411949746
$ws.Range("T4").Value = 0.7906673411949746

# Row 5
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 0.5364909999999999
$ws.Range("H5").Value = 1.609473
$ws.Range("M5").Value = 0.042388
$ws.Range("N5").Value = 0.127164
$ws.Range("O5").Value = 0.0001220879211991998
$ws.Range("P5").Value = 0.0001220879211991998
$ws.Range("Q5").Value = 0.022740780508
$ws.Range("R5").Value = 0.204667024572
$ws.Range("S5").Value = 0.0001220879211991998
$ws.Range("T5").Value = 0.0001220879211991998

# Row 6
$ws.Range("E6").Value = 2
$ws.Range("F6").Value = 0.6666666666666666
$ws.Range("G6").Value = 0.5364909999999999
$ws.Range("H6").Value = 1.609473
$ws.Range("M6").Value = 72.37795533333333
$ws.Range("N6").Value = 217.133866
$ws.Range("O6").Value = 0.2084664081177503
$ws.Range("P6").Value = 0.2084664081177503
$ws.Range("Q6").Value = 38.83012163473533
$ws.Range("R6").Value = 349.471094712618
$ws.Range("S6").Value = 0.2084664081177503
$ws.Range("T6").Value = 0.2084664081177503
